$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 5-7 (ECs sending-cluster rows); new TPM data only keeps
# the MuSCs sending-cluster rows (previously rows 2-4) with updated numbers.
$ws.Range("A5:A7").EntireRow.Delete()

# Row 2: MuSCs / Gdnf / Gfra1 / ECs
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Gdnf"
$ws.Range("C2").Value = "Gfra1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.509729
$ws.Range("H2").Value = 1.529187
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.03257366666666667
$ws.Range("N2").Value = 0.097721
$ws.Range("O2").Value = 0.001227793554179957
$ws.Range("P2").Value = 0.001227793554179957
$ws.Range("Q2").Value = 0.01660374253633333
$ws.Range("R2").Value = 0.149433682827
$ws.Range("S2").Value = 0.001227793554179957
$ws.Range("T2").Value = 0.001227793554179957

# Row 3: MuSCs / Gdnf / Gfra1 / FAPs
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Gdnf"
$ws.Range("C3").Value = "Gfra1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.509729
$ws.Range("H3").Value = 1.529187
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 20.328499
$ws.Range("N3").Value = 60.985497
$ws.Range("O3").Value = 0.7662385783512358
$ws.Range("P3").Value = 0.7662385783512359
$ws.Range("Q3").Value = 10.362025466771
$ws.Range("R3").Value = 93.25822920093898
$ws.Range("S3").Value = 0.7662385783512358
$ws.Range("T3").Value = 0.7662385783512359

# Row 4: MuSCs / Gdnf / Gfra1 / MuSCs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Gdnf"
$ws.Range("C4").Value = "Gfra1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.509729
$ws.Range("H4").Value = 1.529187
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.169174666666667
$ws.Range("N4").Value = 18.507524
$ws.Range("O4").Value = 0.2325336280945842
$ws.Range("P4").Value = 0.2325336280945842
$ws.Range("Q4").Value = 3.144607233665333
$ws.Range("R4").Value = 28.301465102988
$ws.Range("S4").Value = 0.2325336280945842
$ws.Range("T4").Value = 0.2325336280945842
